$wb = $excel.ActiveWorkbook

# "Test Cases" sheet: row 4 (Login_03) Result column D4 changes from PASS to FAIL
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("D4").Value = "FAIL"

# "Test Steps" sheet: H12, H17, H21 Result values change from PASS to FAIL
$wsTestSteps = $wb.Worksheets.Item("Test Steps")
$wsTestSteps.Range("H12").Value = "FAIL"
$wsTestSteps.Range("H17").Value = "FAIL"
$wsTestSteps.Range("H21").Value = "FAIL"
